# Update cryptocurrency price/volume data to reflect the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $value) {
    $range = $ws.Range($cellRef)
    # Force a text number format so numeric-looking strings (e.g. '1.00')
    # are not silently reinterpreted as numbers by Excel, then clear the
    # formatting change again so the cell keeps its original (default) style.
    $range.NumberFormat = '@'
    $range.Value = $value
    $range.ClearFormats()
}


# Row 2
Set-TextCell 'D2' '66.999.12'
Set-TextCell 'E2' '  +0.70%  '

# Row 3
Set-TextCell 'D3' '3.797.04'
Set-TextCell 'E3' '  -1.37%  '

# Row 4
Set-TextCell 'D4' '1.00'
Set-TextCell 'E4' '  +0.11%  '

# Row 5
Set-TextCell 'D5' '436.09'
Set-TextCell 'E5' '  +1.30%  '

# Row 6
Set-TextCell 'D6' '139.74'
Set-TextCell 'E6' '  +6.57%  '

# Row 7
Set-TextCell 'D7' '0.622'
Set-TextCell 'E7' '  +1.87%  '

# Row 8
Set-TextCell 'D8' '1.00'
Set-TextCell 'E8' '  +0.01%  '

# Row 9
Set-TextCell 'D9' '0.737'
Set-TextCell 'E9' '  +1.24%  '

# Row 10
Set-TextCell 'E10' '  -7.51%  '

# Row 11
Set-TextCell 'D11' '0.0000319'
Set-TextCell 'E11' '  -12.10%  '

# Row 12
Set-TextCell 'D12' '42.88'
Set-TextCell 'E12' '  +5.13%  '

# Row 13
Set-TextCell 'D13' '10.44'
Set-TextCell 'E13' '  +4.04%  '

# Row 14
Set-TextCell 'D14' '4.404.48'
Set-TextCell 'E14' '  -1.12%  '

# Row 15
Set-TextCell 'D15' '14.86'
Set-TextCell 'E15' '  -5.98%  '

# Row 16
Set-TextCell 'E16' '  -0.50%  '

# Row 17
Set-TextCell 'D17' '3.812.91'
Set-TextCell 'E17' '  -0.67%  '

# Row 18
Set-TextCell 'D18' '19.95'
Set-TextCell 'E18' '  +2.02%  '

# Row 19
Set-TextCell 'D19' '1.13'
Set-TextCell 'E19' '  +6.45%  '

# Row 20
Set-TextCell 'D20' '67.093.54'
Set-TextCell 'E20' '  +0.35%  '

# Row 21
Set-TextCell 'D21' '419.82'
Set-TextCell 'E21' '  +2.94%  '

# Row 22
Set-TextCell 'D22' '14.69'
Set-TextCell 'E22' '  +1.98%  '

# Row 23
Set-TextCell 'E23' '  +6.80%  '

# Row 24
Set-TextCell 'D24' '86.19'
Set-TextCell 'E24' '  +1.27%  '

# Row 25
Set-TextCell 'D25' '37.30'
Set-TextCell 'E25' '  +0.58%  '

# Row 26
Set-TextCell 'D26' '3.39'
Set-TextCell 'E26' '  +4.56%  '

# Row 27
Set-TextCell 'B27' 'LEO'
Set-TextCell 'C27' 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextCell 'D27' '5.61'
Set-TextCell 'E27' '  -0.90%  '

# Row 28
Set-TextCell 'B28' 'Filecoin'
Set-TextCell 'C28' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextCell 'D28' '9.90'
Set-TextCell 'E28' '  +4.25%  '

# Row 29
Set-TextCell 'B29' 'RenderToken'
Set-TextCell 'C29' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell 'D29' '9.61'
Set-TextCell 'E29' '  +34.55%  '

# Row 30
Set-TextCell 'D30' '741.70'
Set-TextCell 'E30' '  +7.68%  '

# Row 31
Set-TextCell 'E31' '  +10.99%  '

# Row 32
Set-TextCell 'D32' '0.133'
Set-TextCell 'E32' '  +10.02%  '

# Row 33
Set-TextCell 'E33' '  +2.65%  '

# Row 34
Set-TextCell 'D34' '43.71'
Set-TextCell 'E34' '  +12.79%  '

# Row 35
Set-TextCell 'E35' '  +4.77%  '

# Row 36
Set-TextCell 'E36' '  +0.07%  '

# Row 37
Set-TextCell 'D37' '5.59'
Set-TextCell 'E37' '  +22.86%  '

# Row 38
Set-TextCell 'D38' '56.32'
Set-TextCell 'E38' '  +2.02%  '

# Row 39
Set-TextCell 'D39' '0.0480'
Set-TextCell 'E39' '  +5.16%  '

# Row 40
Set-TextCell 'D40' '2.67'
Set-TextCell 'E40' '  +38.90%  '

# Row 41
Set-TextCell 'D41' '2.94'
Set-TextCell 'E41' '  -4.10%  '

# Row 42
Set-TextCell 'D42' '0.0₃0686'
Set-TextCell 'E42' '  -14.78%  '

# Row 43
Set-TextCell 'E43' '  +3.82%  '

# Row 44
Set-TextCell 'B44' 'FirstDigitalUSD'
Set-TextCell 'C44' 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextCell 'D44' '1.00'
Set-TextCell 'E44' '  -0.11%  '

# Row 45
Set-TextCell 'B45' 'ApeXProtocol'
Set-TextCell 'C45' 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextCell 'D45' '3.34'
Set-TextCell 'E45' '  +7.08%  '

# Row 46
Set-TextCell 'D46' '0.330'
Set-TextCell 'E46' '  +12.61%  '

# Row 47
Set-TextCell 'D47' '3.31'
Set-TextCell 'E47' '  +1.16%  '

# Row 49
Set-TextCell 'D49' '2.66'
Set-TextCell 'E49' '  +4.60%  '

# Row 50
Set-TextCell 'D50' '143.10'
Set-TextCell 'E50' '  -3.35%  '

# Row 51
Set-TextCell 'E51' '  +2.17%  '
